# Mental Health Chatbot.pptx - apply commit changes
#
# 1) Reposition the "Blanche Payton" textbox (shape id 24) to its new
#    location on the slide.
# 2) Add a new textbox ("TextBox 3") containing the project's GitHub URL,
#    hyperlinked back to slide 1 (a "link to slide in this presentation"
#    hyperlink, as created when selecting text and inserting a hyperlink
#    to a slide).
#
# Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU)
# by this object model, and the underlying engine truncates when
# converting points -> EMU, so a tiny fraction of a point is added before
# the conversion to land exactly on the target EMU value.

function EmuToPt($emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

# AddTextbox's own EMU conversion rounds rather than truncates, so no
# epsilon is needed there.
function EmuToPtExact($emu) {
    return $emu / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Move the "Blanche Payton" name textbox ---------------------------
$nameBox = $s.Shapes.Item(9)
$nameBox.Left = EmuToPt 591911
$nameBox.Top  = EmuToPt 6211290

# --- 2. Add the new GitHub-link textbox -----------------------------------
$linkBox = $s.Shapes.AddTextbox(
    1,
    (EmuToPtExact 6225432),
    (EmuToPtExact 5934660),
    (EmuToPtExact 4810868),
    (EmuToPtExact 923330)
)

$linkBox.Fill.Visible = $false
$linkBox.TextFrame.WordWrap = $true
$linkBox.TextFrame.AutoSize = 1

$tr = $linkBox.TextFrame.TextRange
$tr.Text = "https://github.com/ksu-hmi/AI-Chatbot-for-Mental-Health-Support-and-Medication-Education"

$runRange = $tr.Characters(1, $tr.Length)
$hyperlink = $runRange.ActionSettings.Item(1).Hyperlink
$hyperlink.Address = ""
$hyperlink.SubAddress = "1"

# Setting AutoSize re-fits the box height to the (single default-sized
# line of) text, so restore the authored height afterwards.
$linkBox.Height = EmuToPt 923330

Write-Host "Applied Mental Health Chatbot.pptx edits"
